# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for "Terminal Hortofrutícola Agro Chillán - Limón"
# before the existing row 641, shifting the previous rows 641-649 down to 645-653.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 641 (pushes old 641..649 down to 645..653)
$ws.Range("A641:A644").EntireRow.Insert()

$data = @(
  @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44656, 16, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "1a amarillo", 100, 21000, 22000, 21500, "`$/malla 16 kilos", "Región de O'Higgins", 1344, 16),
  @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44656, 16, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "1a plateado", 120, 19000, 20000, 19500, "`$/malla 16 kilos", "Región de O'Higgins", 1219, 16),
  @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44656, 16, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "2a amarillo", 100, 19000, 20000, 19500, "`$/malla 16 kilos", "Región de O'Higgins", 1219, 16),
  @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44656, 16, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "2a plateado", 80, 18000, 18000, 18000, "`$/malla 16 kilos", "Región de O'Higgins", 1125, 16)
)

for ($i = 0; $i -lt 4; $i++) {
  $r = 641 + $i
  $row = $data[$i]
  for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}
